# delete dissolved municipalities from options (#100)
# "Beemster" municipality was dissolved and merged into Purmerend, so its
# row is removed from the "area" options sheet.

$wb = $excel.ActiveWorkbook
$wsOutcome = $wb.Worksheets.Item("outcome")
$wsArea = $wb.Worksheets.Item("area")

# Row 57 on the "area" sheet holds the "Beemster" / "Municipality" option.
# Delete the whole row; everything below shifts up automatically (and the
# now-unused "Beemster" shared string is dropped from the workbook).
[void]$wsArea.Rows(57).Delete()

# Restore the selections recorded in the sheets, and leave the "area" sheet
# as the active / selected tab, matching the saved workbook state.
[void]$wsOutcome.Activate()
[void]$wsOutcome.Range("E12").Select()

[void]$wsArea.Activate()
[void]$wsArea.Range("K53").Select()
